# The "numberOfPages" column (column F) was removed from the manifest
# sheet. Select the column first (so the saved selection/active-cell
# mirrors what Excel records after an interactive column delete), then
# delete it - this shifts every column from G onward one slot to the
# left (G->F, H->G, ... S->R) and drops the now-unused "numberOfPages"
# shared string automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns("F").Select() | Out-Null
$ws.Columns("F").EntireColumn.Delete()
